# edit.ps1 -- apply the "Quantum Entanglement" -> "Chemistry" essay rewrite
# (title / author / email / body / summary) plus the TimesNewToman ->
# Times New Roman font-name fix, as described by the target diff.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $old"
    }
}

# ---------------------------------------------------------------------
# 1. Font: TimesNewToman -> Times New Roman, everywhere in the body.
#    (Setting Find.Font.Name doesn't reliably match here, so walk every
#    non-empty paragraph range and set its Font.Name directly -- but
#    stop one character short of the paragraph mark so we don't stamp
#    a spurious rPr onto the pPr of every paragraph.)
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $start = $p.Range.Start
    $end = $p.Range.End - 1
    if ($end -gt $start) {
        $body = $d.Range($start, $end)
        $body.Font.Name = "Times New Roman"
    }
}

# ---------------------------------------------------------------------
# 2. Title
# ---------------------------------------------------------------------
Replace-Text "The Enigmatic Realm of Quantum Entanglement" "Exploring the Marvels of Chemistry: A Journey into the Molecular World"

# ---------------------------------------------------------------------
# 3. Byline / author name
# ---------------------------------------------------------------------
Replace-Text "Emily Jones" "Dr. Lucy Thompson"

# ---------------------------------------------------------------------
# 4. Email address line
# ---------------------------------------------------------------------
Replace-Text "Emily@QuantumStudies" "Lucy"
Replace-Text "edu" "Thompson@eduworld.org"

# ---------------------------------------------------------------------
# 5. Body paragraph sentences (in reading order)
# ---------------------------------------------------------------------
Replace-Text "Within the ethereal expanse of quantum mechanics, a peculiar phenomenon known as entanglement captivates the imagination and challenges the very foundations of our understanding of reality" "Chemistry, the study of matter and its properties, is a fascinating field that unveils the fundamental building blocks of our universe"

Replace-Text " This intricate dance between particles, where the state of one instantaneously influences the state of another, regardless of the distance between them, has profound implications that have ignited fervent debate among physicists, philosophers, and even artists" " In this realm of atoms, molecules, and chemical reactions, we embark on a journey to explore the intricate tapestry of substances that shape our world"

Replace-Text " Enter the enigmatic realm of quantum entanglement, where the boundaries of space and time seem to dissolve, blurring the lines between interconnectedness and individuality" " From the air we breathe to the food we eat, chemistry plays a pivotal role in understanding the phenomena that govern our everyday lives"

Replace-Text "Unveiling the enigmatic nature of entanglement has far-reaching consequences, not just in the abstract realm of theoretical physics, but also in the practical applications that are shaping our technological landscape" "As alchemists of modern times, chemists strive to comprehend the enigmatic forces that bind atoms, unraveling the secrets of their interactions"

Replace-Text " From the development of quantum computers that promise exponential leaps in computational power to the secure communication networks that safeguard our digital interactions, entanglement is poised to revolutionize industries and redefine the way we communicate, compute, and perceive the world around us" " Through careful experimentation and analysis, they unlock the mysteries of chemical reactions, revealing the symphony of energy transformations that drive countless processes in nature"

# Insert the two brand-new runs (" It is in this dance..." + ".") right
# after the sentence/period that now reads "...processes in nature."
# and before the following line break.
$anchor = $d.Content
$anchor.Find.Execute("Through careful experimentation and analysis, they unlock the mysteries of chemical reactions, revealing the symphony of energy transformations that drive countless processes in nature")
$anchor.Collapse(0)
$anchor.MoveEnd(1, 1)
$anchor.Collapse(0)
$anchor.InsertAfter(" It is in this dance of molecules that we discover the exquisite beauty and elegance of the chemical world.")

Replace-Text "As we delve deeper into the mysteries of entanglement, we find ourselves confronted with profound questions that touch upon the very nature of reality" "Furthermore, chemistry has revolutionized various industries, leading to advancements in materials science, pharmaceuticals, and energy production"

Replace-Text " Some posit that entanglement provides evidence for a deeper level of interconnectedness in the universe, hinting at a hidden order yet to be fully grasped" " The synthesis of new materials with tailored properties has paved the way for technological breakthroughs in fields ranging from electronics to aerospace"

Replace-Text " Others contend that it challenges our classical notions of locality and causality, forcing us to reconsider the fundamental assumptions upon which our understanding of the cosmos is built" " Similarly, the development of life-saving drugs and therapies has transformed medicine, improving the quality of life for millions worldwide"

# ---------------------------------------------------------------------
# 6. Summary paragraph
# ---------------------------------------------------------------------
Replace-Text "Quantum entanglement stands as a testament to the mesmerizing strangeness of the quantum realm, a phenomenon that has ignited both scientific exploration and philosophical contemplation" "This essay provided a glimpse into the captivating world of chemistry, showcasing its fundamental importance in understanding the nature of matter, chemical reactions, and their wide-ranging applications"

Replace-Text " Its potential applications hold the promise of transformative technologies, while its deeper implications continue to challenge our understanding of the fundamental nature of reality" " From the intricacies of molecular interactions to the practical implications in various industries, chemistry stands as a testament to the power of scientific inquiry and its transformative impact on society"

# Drop the final sentence + its trailing period entirely -- the summary
# now ends after the run above's period.
$tail = $d.Content
$tail.Find.Execute(" As we unravel the intricacies of entanglement, we are embarking on a journey into the unknown, pushing the boundaries of human knowledge and perhaps, one day, glimpsing the hidden harmonies of the universe")
$tail.MoveEnd(1, 1)
$tail.Text = ""

# ---------------------------------------------------------------------
# 7. Trailing empty paragraph added after the summary, before the
#    section break.
# ---------------------------------------------------------------------
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()
